{"js": "// Change the year in the astromap link: 2018 -> 2022.\n// The paragraph that credits the map author ends with a sentence whose\n// runs are collapsed into a single plain run (preceded by one leftover\n// empty run), matching how Word normalizes a Find&Replace-style edit\n// across many runs.\n\nconst OLD_SENTENCE =\n  \"T\u00e4m\u00e4n oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).\";\nconst NEW_SENTENCE =\n  \"T\u00e4m\u00e4n oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"T\u00e4m\u00e4n oppaan kartat piirsi\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the paragraph crediting the map author.\");\n}\n\n// Build a minimal OOXML fragment for a paragraph containing a leading\n// empty run plus one plain run with the updated sentence. Inserting it\n// \"before\" the paragraph's own Start keeps the paragraph's own\n// properties (pPr: borders/centering/etc.) untouched, and only adds the\n// new runs ahead of the old ones.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r/><w:r><w:t>' +\n  NEW_SENTENCE +\n  '</w:t></w:r></w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.getRange(\"Start\").insertOoxml(ooxml, Word.InsertLocation.before);\nawait context.sync();\n\n// Remove the old (now duplicated) sentence text that still trails the\n// newly-inserted runs, leaving only the new runs behind.\nconst oldResults = target.search(OLD_SENTENCE, { matchCase: true });\noldResults.load(\"items\");\nawait context.sync();\n\nif (oldResults.items.length === 0) {\n  throw new Error(\"Could not find the old sentence text to remove.\");\n}\n\noldResults.items[0].delete();\nawait context.sync();\n", "ps1": "# Change the year in the astromap link: 2018 -> 2022.\n# The paragraph that credits the map author ends with a sentence whose\n# runs are collapsed into a single plain run (preceded by one leftover\n# empty run), matching how Word normalizes a Find&Replace-style edit\n# across many runs.\n\n$d = $word.ActiveDocument\n\n$oldSentence = \"T\u00e4m\u00e4n oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).\"\n$newSentence = \"T\u00e4m\u00e4n oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*T\u00e4m\u00e4n oppaan kartat piirsi*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the paragraph crediting the map author.\"\n}\n\n$r = $target.Range\n$oldStart = $r.Start\n$oldEnd = $r.End - 1   # exclude the trailing paragraph mark\n\n# Insert the new (plain-formatted) sentence directly before the old\n# content, then drop an empty marker run at the original start so the\n# final paragraph begins with one leftover empty run, matching the\n# target structure.\n$r.InsertBefore($newSentence)\n$marker = $d.Range($oldStart, $oldStart)\n$marker.InsertBefore(\"\")\n\n# Now delete the old sentence text, which got pushed past the newly\n# inserted text.\n$insertedLen = $newSentence.Length\n$newOldStart = $oldStart + $insertedLen\n$newOldEnd = $oldEnd + $insertedLen\n$delRange = $d.Range($newOldStart, $newOldEnd)\n$delRange.Delete()\n"}
